$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 7.783173581191642
$ws.Cells.Item(2, 4).Value = 7.839920480229148
$ws.Cells.Item(2, 5).Value = 12.02209393141306
$ws.Cells.Item(2, 6).Value = 40.04974693950228
$ws.Cells.Item(2, 7).Value = 50.40278016977956
$ws.Cells.Item(2, 8).Value = 18.53676769311976
$ws.Cells.Item(2, 10).Value = 9.510999597123126
$ws.Cells.Item(2, 13).Value = 23.80390300648689
$ws.Cells.Item(2, 14).Value = 17.81153797085744
$ws.Cells.Item(3, 2).Value = 7.709723325812419
$ws.Cells.Item(3, 4).Value = 7.808524432377236
$ws.Cells.Item(3, 5).Value = 12.05226414689429
$ws.Cells.Item(3, 6).Value = 39.67740257645623
$ws.Cells.Item(3, 7).Value = 49.37724731733471
$ws.Cells.Item(3, 8).Value = 18.45225351759321
$ws.Cells.Item(3, 10).Value = 9.548714958702895
$ws.Cells.Item(3, 13).Value = 23.18563396420499
$ws.Cells.Item(3, 14).Value = 17.7738682296755
$ws.Cells.Item(4, 2).Value = 7.666071942561167
$ws.Cells.Item(4, 4).Value = 7.788941521657875
$ws.Cells.Item(4, 5).Value = 12.07186712546299
$ws.Cells.Item(4, 6).Value = 39.46191230231511
$ws.Cells.Item(4, 7).Value = 48.75812006814865
$ws.Cells.Item(4, 8).Value = 18.40603214316965
$ws.Cells.Item(4, 10).Value = 9.572981056780449
$ws.Cells.Item(4, 13).Value = 22.80086583612374
$ws.Cells.Item(4, 14).Value = 17.75325082035187
$ws.Cells.Item(5, 2).Value = 7.648667800845041
$ws.Cells.Item(5, 4).Value = 7.78088504315551
$ws.Cells.Item(5, 5).Value = 12.08012739049211
$ws.Cells.Item(5, 6).Value = 39.37747959496287
$ws.Cells.Item(5, 7).Value = 48.50886067323514
$ws.Cells.Item(5, 8).Value = 18.38863130756727
$ws.Cells.Item(5, 10).Value = 9.583149466956103
$ws.Cells.Item(5, 13).Value = 22.64299723067514
$ws.Cells.Item(5, 14).Value = 17.74548818592557
$ws.Cells.Item(6, 2).Value = 7.645801634713722
$ws.Cells.Item(6, 4).Value = 7.779542679915674
$ws.Cells.Item(6, 5).Value = 12.08151544438165
$ws.Cells.Item(6, 6).Value = 39.36366599527291
$ws.Cells.Item(6, 7).Value = 48.4676661714293
$ws.Cells.Item(6, 8).Value = 18.38582878053411
$ws.Cells.Item(6, 10).Value = 9.584854850465447
$ws.Cells.Item(6, 13).Value = 22.61672521111601
$ws.Cells.Item(6, 14).Value = 17.74423801883558
$ws.Cells.Item(7, 2).Value = 7.665835642623647
$ws.Cells.Item(7, 4).Value = 7.788833177310882
$ws.Cells.Item(7, 5).Value = 12.07197742442069
$ws.Cells.Item(7, 6).Value = 39.46075982330647
$ws.Cells.Item(7, 7).Value = 48.75474565182966
$ws.Cells.Item(7, 8).Value = 18.40579165015515
$ws.Cells.Item(7, 10).Value = 9.57311705738991
$ws.Cells.Item(7, 13).Value = 22.79874080528145
$ws.Cells.Item(7, 14).Value = 17.75314353334516
$ws.Cells.Item(8, 2).Value = 7.757559854521983
$ws.Cells.Item(8, 4).Value = 7.829158889978405
$ws.Cells.Item(8, 5).Value = 12.03227325388373
$ws.Cells.Item(8, 6).Value = 39.91868748203282
$ws.Cells.Item(8, 7).Value = 50.04721809505018
$ws.Cells.Item(8, 8).Value = 18.50645621963617
$ws.Cells.Item(8, 10).Value = 9.523774334834009
$ws.Cells.Item(8, 13).Value = 23.59192201790287
$ws.Cells.Item(8, 14).Value = 17.79803268088581
$ws.Cells.Item(9, 2).Value = 7.947996968698414
$ws.Cells.Item(9, 4).Value = 7.905789688955638
$ws.Cells.Item(9, 5).Value = 11.96293626252313
$ws.Cells.Item(9, 6).Value = 40.91699769271266
$ws.Cells.Item(9, 7).Value = 52.64853613396208
$ws.Cells.Item(9, 8).Value = 18.74836345054039
$ws.Cells.Item(9, 10).Value = 9.435765175567516
$ws.Cells.Item(9, 13).Value = 25.09735941128893
$ws.Cells.Item(9, 14).Value = 17.90564832874191
$ws.Cells.Item(10, 2).Value = 8.093089462734838
$ws.Cells.Item(10, 4).Value = 7.960570364428005
$ws.Cells.Item(10, 5).Value = 11.91714341331563
$ws.Cells.Item(10, 6).Value = 41.70595765447587
$ws.Cells.Item(10, 7).Value = 54.57809004680612
$ws.Cells.Item(10, 8).Value = 18.95242519426461
$ws.Cells.Item(10, 10).Value = 9.376375902592486
$ws.Cells.Item(10, 13).Value = 26.16163077427799
$ws.Cells.Item(10, 14).Value = 17.99617406619567
$ws.Cells.Item(11, 2).Value = 8.159945912094404
$ws.Cells.Item(11, 4).Value = 7.985155530807855
$ws.Cells.Item(11, 5).Value = 11.89741919107877
$ws.Cells.Item(11, 6).Value = 42.07561145731312
$ws.Cells.Item(11, 7).Value = 55.45535116030212
$ws.Cells.Item(11, 8).Value = 19.05076352771231
$ws.Cells.Item(11, 10).Value = 9.35048895580565
$ws.Cells.Item(11, 13).Value = 26.63461367131775
$ws.Cells.Item(11, 14).Value = 18.03973531775564
$ws.Cells.Item(12, 2).Value = 8.185362316328961
$ws.Cells.Item(12, 4).Value = 7.994416528965517
$ws.Cells.Item(12, 5).Value = 11.89010859550982
$ws.Cells.Item(12, 6).Value = 42.21701967869757
$ws.Cells.Item(12, 7).Value = 55.78711757215908
$ws.Cells.Item(12, 8).Value = 19.08877315364387
$ws.Cells.Item(12, 10).Value = 9.340847622497263
$ws.Cells.Item(12, 13).Value = 26.81195838864105
$ws.Cells.Item(12, 14).Value = 18.05656330799093
$ws.Cells.Item(13, 2).Value = 8.179884392274969
$ws.Cells.Item(13, 4).Value = 7.992424192411633
$ws.Cells.Item(13, 5).Value = 11.89167602263029
$ws.Cells.Item(13, 6).Value = 42.18650326325506
$ws.Cells.Item(13, 7).Value = 55.71569056469757
$ws.Cells.Item(13, 8).Value = 19.08055316459813
$ws.Cells.Item(13, 10).Value = 9.342916888577674
$ws.Cells.Item(13, 13).Value = 26.77384485046684
$ws.Cells.Item(13, 14).Value = 18.05292448376005
$ws.Cells.Item(14, 2).Value = 8.162035077839201
$ws.Cells.Item(14, 4).Value = 7.985918422290614
$ws.Cells.Item(14, 5).Value = 11.89681457044601
$ws.Cells.Item(14, 6).Value = 42.08721725404867
$ws.Cells.Item(14, 7).Value = 55.48265647008098
$ws.Cells.Item(14, 8).Value = 19.05387527749317
$ws.Cells.Item(14, 10).Value = 9.349692526319062
$ws.Cells.Item(14, 13).Value = 26.6492401343712
$ws.Cells.Item(14, 14).Value = 18.04111316570079
$ws.Cells.Item(15, 2).Value = 8.151114100256091
$ws.Cells.Item(15, 4).Value = 7.98192706432862
$ws.Cells.Item(15, 5).Value = 11.89998270601171
$ws.Cells.Item(15, 6).Value = 42.02658418094016
$ws.Cells.Item(15, 7).Value = 55.33984952234515
$ws.Cells.Item(15, 8).Value = 19.03763405346398
$ws.Cells.Item(15, 10).Value = 9.353863803387682
$ws.Cells.Item(15, 13).Value = 26.57268199283256
$ws.Cells.Item(15, 14).Value = 18.03392136283404
$ws.Cells.Item(16, 2).Value = 8.088735404343128
$ws.Cells.Item(16, 4).Value = 7.958956866651463
$ws.Cells.Item(16, 5).Value = 11.91845465935043
$ws.Cells.Item(16, 6).Value = 41.68200626025839
$ws.Cells.Item(16, 7).Value = 54.52071758001351
$ws.Cells.Item(16, 8).Value = 18.94610764344904
$ws.Cells.Item(16, 10).Value = 9.378090324630897
$ws.Cells.Item(16, 13).Value = 26.13048176587063
$ws.Cells.Item(16, 14).Value = 17.99337428922453
$ws.Cells.Item(17, 2).Value = 8.05066995302467
$ws.Cells.Item(17, 4).Value = 7.944779123016375
$ws.Cells.Item(17, 5).Value = 11.93006969957034
$ws.Cells.Item(17, 6).Value = 41.47328980227752
$ws.Cells.Item(17, 7).Value = 54.0178153873763
$ws.Cells.Item(17, 8).Value = 18.89135573290799
$ws.Cells.Item(17, 10).Value = 9.393241145583488
$ws.Cells.Item(17, 13).Value = 25.85623048672781
$ws.Cells.Item(17, 14).Value = 17.96910266754167
$ws.Cells.Item(18, 2).Value = 8.02885740993784
$ws.Cells.Item(18, 4).Value = 7.936592986959756
$ws.Cells.Item(18, 5).Value = 11.93685461589007
$ws.Cells.Item(18, 6).Value = 41.35426089815189
$ws.Cells.Item(18, 7).Value = 53.72853256375436
$ws.Cells.Item(18, 8).Value = 18.86038385921725
$ws.Cells.Item(18, 10).Value = 9.402061869581514
$ws.Cells.Item(18, 13).Value = 25.69744661829358
$ws.Cells.Item(18, 14).Value = 17.9553669716024
$ws.Cells.Item(19, 2).Value = 8.021486802830676
$ws.Cells.Item(19, 4).Value = 7.933815918359222
$ws.Cells.Item(19, 5).Value = 11.93916979792714
$ws.Cells.Item(19, 6).Value = 41.3141383620619
$ws.Cells.Item(19, 7).Value = 53.63059306480424
$ws.Cells.Item(19, 8).Value = 18.84998724063735
$ws.Cells.Item(19, 10).Value = 9.405066713218911
$ws.Cells.Item(19, 13).Value = 25.64351127411796
$ws.Cells.Item(19, 14).Value = 17.95075519398397
$ws.Cells.Item(20, 2).Value = 8.054713799504805
$ws.Cells.Item(20, 4).Value = 7.946291630561547
$ws.Cells.Item(20, 5).Value = 11.92882247386104
$ws.Cells.Item(20, 6).Value = 41.49540337842407
$ws.Cells.Item(20, 7).Value = 54.07135562935063
$ws.Cells.Item(20, 8).Value = 18.8971305001193
$ws.Cells.Item(20, 10).Value = 9.391617312453562
$ws.Cells.Item(20, 13).Value = 25.88553402035647
$ws.Cells.Item(20, 14).Value = 17.97166323937482
$ws.Cells.Item(21, 2).Value = 8.167275344240011
$ws.Cells.Item(21, 4).Value = 7.987830656795795
$ws.Cells.Item(21, 5).Value = 11.89530095700767
$ws.Cells.Item(21, 6).Value = 42.11634214766852
$ws.Cells.Item(21, 7).Value = 55.55111877935561
$ws.Cells.Item(21, 8).Value = 19.06169047156384
$ws.Cells.Item(21, 10).Value = 9.347697981859717
$ws.Cells.Item(21, 13).Value = 26.68588862163758
$ws.Cells.Item(21, 14).Value = 18.04457349946859
$ws.Cells.Item(22, 2).Value = 8.241409004326362
$ws.Cells.Item(22, 4).Value = 8.014693617723047
$ws.Cells.Item(22, 5).Value = 11.87431648801574
$ws.Cells.Item(22, 6).Value = 42.53043086186955
$ws.Cells.Item(22, 7).Value = 56.51554959862131
$ws.Cells.Item(22, 8).Value = 19.17372310013478
$ws.Cells.Item(22, 10).Value = 9.319935025806075
$ws.Cells.Item(22, 13).Value = 27.19862286708542
$ws.Cells.Item(22, 14).Value = 18.09415653301901
$ws.Cells.Item(23, 2).Value = 8.20179816405664
$ws.Cells.Item(23, 4).Value = 8.000382684783066
$ws.Cells.Item(23, 5).Value = 11.88543198870351
$ws.Cells.Item(23, 6).Value = 42.30870678272957
$ws.Cells.Item(23, 7).Value = 56.00117438867333
$ws.Cells.Item(23, 8).Value = 19.11352640871217
$ws.Cells.Item(23, 10).Value = 9.334666856444269
$ws.Cells.Item(23, 13).Value = 26.92596230216598
$ws.Cells.Item(23, 14).Value = 18.0675197576975
$ws.Cells.Item(24, 2).Value = 8.052885351122358
$ws.Cells.Item(24, 4).Value = 7.945607935368393
$ws.Cells.Item(24, 5).Value = 11.92938601041592
$ws.Cells.Item(24, 6).Value = 41.48540281608688
$ws.Cells.Item(24, 7).Value = 54.04715054941344
$ws.Cells.Item(24, 8).Value = 18.89451815103929
$ws.Cells.Item(24, 10).Value = 9.39235110373386
$ws.Cells.Item(24, 13).Value = 25.87228934782278
$ws.Cells.Item(24, 14).Value = 17.97050492387823
$ws.Cells.Item(25, 2).Value = 7.895476574444182
$ws.Cells.Item(25, 4).Value = 7.885323887897025
$ws.Cells.Item(25, 5).Value = 11.98078617416179
$ws.Cells.Item(25, 6).Value = 40.63674925523083
$ws.Cells.Item(25, 7).Value = 51.94000216199047
$ws.Cells.Item(25, 8).Value = 18.67823080431054
$ws.Cells.Item(25, 10).Value = 9.458643701075568
$ws.Cells.Item(25, 13).Value = 24.6966706714888
$ws.Cells.Item(25, 14).Value = 17.87448656540692
